# Update market_objects sheet: refresh smoothed volsurface / q^probability data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows whose MarketObjects lists changed
$ws.Range("B93").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"
$ws.Range("B94").Value = "['BTCUSD.SPOT']"
$ws.Range("B95").Value = "['BTCUSD.SPOT']"

# Append new rows 96-100
$newRows = @(
    @{ Row = 96;  Date = "2025-09-15"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']" },
    @{ Row = 97;  Date = "2025-09-16"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']" },
    @{ Row = 98;  Date = "2025-09-17"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']" },
    @{ Row = 99;  Date = "2025-09-18"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 100; Date = "2025-09-19"; Objects = "['USD.SOFR.CSA_USD']" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.Objects
}
